# Generate Report for Handback
# The file "8e59d3a0-f493-47bb-8bdd-835f569b2adb.md" (zh-cn and de-de rows 3)
# has moved from "Ready for handoff" to "Handed back: in sync with en-US",
# and its handback timestamps are refreshed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update the status shown for zh-cn / de-de columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn detail sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H2").Value = "2016-03-18 03:27:50"
$wsZhCn.Range("H3").Value = "2016-03-18 03:27:50"

# --- de-de detail sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H2").Value = "2016-03-18 03:28:03"
$wsDeDe.Range("H3").Value = "2016-03-18 03:28:03"
